# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G holds the "K" series; recompute and overwrite the values for rows 2-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 3
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 3
    12 = 2
    13 = 3
    14 = 2
}

foreach ($row in $kValues.Keys | Sort-Object) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
